$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are text that look numeric (e.g. "3.08", "51.658.03").
# Excel would auto-convert such text to a real number on assignment, so we
# force the whole Price column to Text format first, assign the values, then
# restore the "Normal" style so no stray number-format style id is left on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.658.03"
$ws.Range("D3").Value = "3.036.31"
$ws.Range("D5").Value = "387.43"
$ws.Range("D6").Value = "102.83"
$ws.Range("D10").Value = "36.73"
$ws.Range("D13").Value = "3.510.84"
$ws.Range("D14").Value = "18.55"
$ws.Range("D16").Value = "3.030.42"
$ws.Range("D17").Value = "0.974"
$ws.Range("D18").Value = "10.67"
$ws.Range("D19").Value = "51.657.22"
$ws.Range("D20").Value = "3.08"
$ws.Range("D21").Value = "12.51"
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("D23").Value = "69.99"
$ws.Range("D24").Value = "267.52"
$ws.Range("D25").Value = "3.19"
$ws.Range("D26").Value = "8.46"
$ws.Range("D27").Value = "7.45"
$ws.Range("D30").Value = "26.26"
$ws.Range("D34").Value = "34.11"
$ws.Range("D35").Value = "50.54"
$ws.Range("D36").Value = "0.0448"
$ws.Range("D38").Value = "3.33"
$ws.Range("D39").Value = "0.296"
$ws.Range("D43").Value = "127.07"
$ws.Range("D46").Value = "21.68"
$ws.Range("D49").Value = "2.026.52"
$ws.Range("D50").Value = "3.333.57"
$ws.Range("D51").Value = "0.517"

$ws.Range("D2:D51").Style = "Normal"

# Column E (Volume/1h change) values are padded with spaces and a "%" sign,
# so Excel keeps them as plain text automatically.
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("E18").Value = "  -12.01%  "
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  -4.10%  "
$ws.Range("E26").Value = "  +6.73%  "
$ws.Range("E27").Value = "  +5.46%  "
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("E39").Value = "  +8.58%  "
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  +3.85%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("E51").Value = "  +6.38%  "
